$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0.664
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.0054
$ws.Range("F3").Value = 0.6042
$ws.Range("G3").Value = 0.001
$ws.Range("H3").Value = 0.0014
$ws.Range("I3").Value = 0.0004
$ws.Range("J3").Value = 0.0016
$ws.Range("D4").Value = 0.9694
$ws.Range("E4").Value = 0.1928
$ws.Range("F4").Value = 0.5374
$ws.Range("G4").Value = 0.1396
$ws.Range("H4").Value = 0.1534
$ws.Range("I4").Value = 0.1166
$ws.Range("J4").Value = 0.161
$ws.Range("D5").Value = 0.0318
$ws.Range("E5").Value = 0.8094
$ws.Range("F5").Value = 0.4798
$ws.Range("G5").Value = 0.8648
$ws.Range("H5").Value = 0.8488
$ws.Range("I5").Value = 0.8816
$ws.Range("J5").Value = 0.8396
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.9948
$ws.Range("F6").Value = 0.4106
$ws.Range("G6").Value = 0.9992
$ws.Range("H6").Value = 0.9984
$ws.Range("I6").Value = 0.9998
$ws.Range("J6").Value = 0.998
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.351
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.2916
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.2486
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.2038
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("D11").Value = 0
$ws.Range("F11").Value = 0.1652
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("D12").Value = 0
$ws.Range("F12").Value = 0.1352
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("D13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("I13").Value = 0.6228
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.9866
$ws.Range("G14").Value = 0.9928
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 0.5614
$ws.Range("J14").Value = 1
$ws.Range("D15").Value = 0.054
$ws.Range("E15").Value = 0.9314
$ws.Range("F15").Value = 0.7396
$ws.Range("G15").Value = 0.7392
$ws.Range("H15").Value = 0.959
$ws.Range("I15").Value = 0.5216
$ws.Range("J15").Value = 0.953
$ws.Range("D16").Value = 0.9164
$ws.Range("E16").Value = 0.12
$ws.Range("F16").Value = 0.3152
$ws.Range("G16").Value = 0.2912
$ws.Range("H16").Value = 0.0806
$ws.Range("I16").Value = 0.4844
$ws.Range("J16").Value = 0.093
$ws.Range("D17").Value = 0.9998
$ws.Range("E17").Value = 0.002
$ws.Range("F17").Value = 0.1166
$ws.Range("G17").Value = 0.0932
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0.4714
$ws.Range("J17").Value = 0.0008
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0.0324
$ws.Range("G18").Value = 0.0306
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0.4438
$ws.Range("J18").Value = 0
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0.015
$ws.Range("G19").Value = 0.0066
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0.436
$ws.Range("J19").Value = 0
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0.0066
$ws.Range("G20").Value = 0.0052
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0.418
$ws.Range("J20").Value = 0
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0.0028
$ws.Range("G21").Value = 0.0018
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0.4118
$ws.Range("J21").Value = 0
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0.0014
$ws.Range("G22").Value = 0.0012
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0.4076
$ws.Range("J22").Value = 0
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0.0006
$ws.Range("G23").Value = 0.0004
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0.3916
$ws.Range("J23").Value = 0
$ws.Range("D24").Value = 0.0542
$ws.Range("E24").Value = 0.0358
$ws.Range("F24").Value = 0.0056
$ws.Range("G24").Value = 0.9802
$ws.Range("H24").Value = 0.0124
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 0.0276
$ws.Range("D25").Value = 0.1484
$ws.Range("E25").Value = 0.119
$ws.Range("F25").Value = 0.0456
$ws.Range("G25").Value = 0.9116
$ws.Range("H25").Value = 0.079
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 0.1112
$ws.Range("D26").Value = 0.3454
$ws.Range("E26").Value = 0.334
$ws.Range("F26").Value = 0.2738
$ws.Range("G26").Value = 0.6912
$ws.Range("H26").Value = 0.3058
$ws.Range("I26").Value = 0.8802
$ws.Range("J26").Value = 0.3306
$ws.Range("D27").Value = 0.6526
$ws.Range("E27").Value = 0.6684
$ws.Range("F27").Value = 0.7344
$ws.Range("G27").Value = 0.3118
$ws.Range("H27").Value = 0.7106
$ws.Range("I27").Value = 0.1264
$ws.Range("J27").Value = 0.6824
$ws.Range("D28").Value = 0.913
$ws.Range("E28").Value = 0.9252
$ws.Range("F28").Value = 0.98
$ws.Range("G28").Value = 0.0504
$ws.Range("H28").Value = 0.9592
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0.9376
$ws.Range("D29").Value = 0.9938
$ws.Range("E29").Value = 0.9958
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 0.0016
$ws.Range("H29").Value = 0.9994
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0.9976
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0.9998
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 1
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 1
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 1
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = 1
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 1
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1
$ws.Range("D33").Value = 1
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 1
$ws.Range("I33").Value = 0
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 0
$ws.Range("I34").Value = 0
